$d = $word.ActiveDocument

# Helper: find a short unique anchor string that ends with the single
# character we want to change, then replace just that last character in
# place (this preserves the surrounding run boundaries/formatting exactly,
# instead of letting Find/Replace merge runs together).
function Replace-LastChar($anchor, $newChar) {
    $rng = $d.Content
    $found = $rng.Find.Execute($anchor, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Anchor not found: $anchor"
    }
    $pos = $rng.End - 1
    $single = $d.Range($pos, $rng.End)
    $single.Text = $newChar
}

# 1) "toutes les ioinctures" -> "toutes les joinctures"
Replace-LastChar "toutes les i" "j"

# 2) "dencocher les ioinctures" -> "dencocher les joinctures"
Replace-LastChar "dencocher les i" "j"

# 3) "Cest pour iecter" -> "Cest pour jecter"
Replace-LastChar "Cest pour i" "j"

# 4) "<tmp>iours" -> "<tmp>jours"
Replace-LastChar "tmp>i" "j"

# 5) " adiouste y aussy du " -> " adjouste y aussy du " (merges 3 runs into 1)
$d.Content.Find.Execute(" adiouste y aussy du ", $true, $false, $false, $false, $false, $true, 1, $false, " adjouste y aussy du ", 2) | Out-Null

# 6) "pas bien ioinct" -> "pas bien joinct"
Replace-LastChar "pas bien i" "j"

# 7) " escailles celle se font" -> " escailles elles se font"
$d.Content.Find.Execute("escailles celle se font", $true, $false, $false, $false, $false, $true, 1, $false, "escailles elles se font", 2) | Out-Null

# 8) " reioindre" -> " rejoindre"
Replace-LastChar " rei" "j"
